$d = $word.ActiveDocument
Write-Host "test"
